$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3873.0557
$ws.Range("J17").Value = 4012.647
$ws.Range("L17").Value = 12037.941
$ws.Range("N17").Value = -12373.941
$ws.Range("H18").Value = 2459.8
$ws.Range("I18").Value = 2459.8
$ws.Range("K18").Value = 2459.8
$ws.Range("M18").Value = -2175.8
$ws.Range("H20").Value = 3399
$ws.Range("I20").Value = 3399
$ws.Range("K20").Value = 3399
$ws.Range("M20").Value = -3169
$ws.Range("H28").Value = 1788.375
$ws.Range("I28").Value = 1788.375
$ws.Range("K28").Value = 1788.375
$ws.Range("M28").Value = -1303.375
$ws.Range("H35").Value = 3399
$ws.Range("I35").Value = 3399
$ws.Range("K35").Value = 3399
$ws.Range("M35").Value = -3020
$ws.Range("H42").Value = 3169.3076
$ws.Range("I42").Value = 1075
$ws.Range("J42").Value = 6520.2
$ws.Range("K42").Value = 3225
$ws.Range("L42").Value = 19560.6
$ws.Range("M42").Value = -2995
$ws.Range("N42").Value = -20020.6
$ws.Range("H58").Value = 684.5
$ws.Range("J58").Value = 2500
$ws.Range("L58").Value = 7500
$ws.Range("N58").Value = -7800
$ws.Range("H62").Value = 4989
$ws.Range("J62").Value = 5831.6665
$ws.Range("L62").Value = 5831.6665
$ws.Range("N62").Value = -7079.6665
$ws.Range("H65").Value = 4989
$ws.Range("J65").Value = 5831.6665
$ws.Range("L65").Value = 29158.3325
$ws.Range("N65").Value = -35398.3325
$ws.Range("H92").Value = 247
$ws.Range("I92").Value = 241.22223
$ws.Range("K92").Value = 241.22223
$ws.Range("M92").Value = 1006.77777
$ws.Range("H96").Value = 7609.5713
$ws.Range("I96").Value = 9628.182000000001
$ws.Range("J96").Value = 208
$ws.Range("K96").Value = 28884.546
$ws.Range("L96").Value = 624
$ws.Range("M96").Value = -27511.546
$ws.Range("N96").Value = -3370
$ws.Range("H112").Value = 3866
$ws.Range("J112").Value = 3999
$ws.Range("L112").Value = 11997
$ws.Range("N112").Value = -14213
$ws.Range("H132").Value = 1553
$ws.Range("I132").Value = 1553
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4659
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2129
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 771.6
$ws.Range("I135").Value = 771.6
$ws.Range("K135").Value = 6944.400000000001
$ws.Range("M135").Value = -4409.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 1735.75
$ws.Range("I31").Value = 1735.75
$ws.Range("K31").Value = 1735.75
$ws.Range("M31").Value = -1441.75
$ws.Range("H32").Value = 4569.1816
$ws.Range("I32").Value = 4569.1816
$ws.Range("K32").Value = 4569.1816
$ws.Range("M32").Value = -4282.1816
$ws.Range("H61").Value = 3597.5186
$ws.Range("I61").Value = 1569
$ws.Range("K61").Value = 1569
$ws.Range("M61").Value = -1357
$ws.Range("H74").Value = 1960
$ws.Range("I74").Value = 2087.3333
$ws.Range("K74").Value = 2087.3333
$ws.Range("M74").Value = -1213.3333
$ws.Range("H76").Value = 40855.4
$ws.Range("J76").Value = 40855.4
$ws.Range("L76").Value = 40855.4
$ws.Range("N76").Value = -41531.4
$ws.Range("H77").Value = 1960
$ws.Range("I77").Value = 2087.3333
$ws.Range("K77").Value = 10436.6665
$ws.Range("M77").Value = -6068.666499999999
$ws.Range("H79").Value = 40855.4
$ws.Range("J79").Value = 40855.4
$ws.Range("L79").Value = 40855.4
$ws.Range("N79").Value = -43195.4
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1874.36
$ws.Range("I132").Value = 1785.7916
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5357.3748
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2827.3748
$ws.Range("N132").Value = -17060
$ws.Range("H135").Value = 534500
$ws.Range("J135").Value = 534500
$ws.Range("L135").Value = 534500
$ws.Range("N135").Value = -544640
$ws.Range("H136").Value = 3597.5186
$ws.Range("I136").Value = 1569
$ws.Range("K136").Value = 4707
$ws.Range("M136").Value = -2157

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 688.44446
$ws.Range("I5").Value = 662.125
$ws.Range("K5").Value = 662.125
$ws.Range("M5").Value = -549.125
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H94").Value = 1756.75
$ws.Range("I94").Value = 1694
$ws.Range("J94").Value = 1945
$ws.Range("K94").Value = 1694
$ws.Range("L94").Value = 1945
$ws.Range("M94").Value = -1243
$ws.Range("N94").Value = -2847
$ws.Range("H99").Value = 2457.6
$ws.Range("I99").Value = 2064.111
$ws.Range("K99").Value = 2064.111
$ws.Range("M99").Value = -566.1109999999999
$ws.Range("H105").Value = 2564.7778
$ws.Range("I105").Value = 2661
$ws.Range("K105").Value = 2661
$ws.Range("M105").Value = -914

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 919
$ws.Range("I16").Value = 707
$ws.Range("K16").Value = 707
$ws.Range("M16").Value = -420
$ws.Range("H17").Value = 308
$ws.Range("I17").Value = 308
$ws.Range("K17").Value = 308
$ws.Range("M17").Value = -134
$ws.Range("H22").Value = 799.5
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 600
$ws.Range("M22").Value = -250
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 2040.5555
$ws.Range("I31").Value = 1682.4546
$ws.Range("J31").Value = 2603.2856
$ws.Range("K31").Value = 1682.4546
$ws.Range("L31").Value = 2603.2856
$ws.Range("M31").Value = -1387.4546
$ws.Range("N31").Value = -3193.2856
$ws.Range("H34").Value = 2040.5555
$ws.Range("I34").Value = 1682.4546
$ws.Range("J34").Value = 2603.2856
$ws.Range("K34").Value = 1682.4546
$ws.Range("L34").Value = 2603.2856
$ws.Range("M34").Value = -1480.4546
$ws.Range("N34").Value = -3007.2856
$ws.Range("H35").Value = 1544.8
$ws.Range("I35").Value = 1544.8
$ws.Range("K35").Value = 1544.8
$ws.Range("M35").Value = -1250.8
$ws.Range("H58").Value = 1329.5769
$ws.Range("I58").Value = 1329.5769
$ws.Range("K58").Value = 1329.5769
$ws.Range("M58").Value = -1126.5769
$ws.Range("H99").Value = 6192.2
$ws.Range("I99").Value = 4987
$ws.Range("J99").Value = 6995.6665
$ws.Range("K99").Value = 4987
$ws.Range("L99").Value = 6995.6665
$ws.Range("M99").Value = -3489
$ws.Range("N99").Value = -9991.666499999999
$ws.Range("H100").Value = 123920
$ws.Range("J100").Value = 123920
$ws.Range("L100").Value = 123920
$ws.Range("N100").Value = -126084
$ws.Range("H113").Value = 919
$ws.Range("I113").Value = 707
$ws.Range("K113").Value = 707
$ws.Range("M113").Value = 1463
$ws.Range("H122").Value = 1908.8
$ws.Range("I122").Value = 515.6667
$ws.Range("J122").Value = 3998.5
$ws.Range("K122").Value = 1547.0001
$ws.Range("L122").Value = 11995.5
$ws.Range("M122").Value = 902.9999
$ws.Range("N122").Value = -16895.5
$ws.Range("H126").Value = 6192.2
$ws.Range("I126").Value = 4987
$ws.Range("J126").Value = 6995.6665
$ws.Range("K126").Value = 14961
$ws.Range("L126").Value = 20986.9995
$ws.Range("M126").Value = -12491
$ws.Range("N126").Value = -25926.9995
$ws.Range("H133").Value = 36080
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 1329.5769
$ws.Range("I136").Value = 1329.5769
$ws.Range("K136").Value = 3988.7307
$ws.Range("M136").Value = -1438.7307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2872.875
$ws.Range("I34").Value = 296
$ws.Range("K34").Value = 888
$ws.Range("M34").Value = -804
$ws.Range("H46").Value = 19051418
$ws.Range("I46").Value = 104774820
$ws.Range("J46").Value = 1775.3334
$ws.Range("K46").Value = 314324460
$ws.Range("L46").Value = 5326.0002
$ws.Range("M46").Value = -314324369
$ws.Range("N46").Value = -5508.0002
$ws.Range("H109").Value = 1803.8667
$ws.Range("I109").Value = 1843.4286
$ws.Range("K109").Value = 5530.2858
$ws.Range("M109").Value = -4490.2858
$ws.Range("H114").Value = 2281.6155
$ws.Range("I114").Value = 2531.4
$ws.Range("J114").Value = 2125.5
$ws.Range("K114").Value = 7594.200000000001
$ws.Range("L114").Value = 6376.5
$ws.Range("M114").Value = -4340.200000000001
$ws.Range("N114").Value = -12884.5
$ws.Range("H122").Value = 100
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 529153.4399999999
$ws.Range("J131").Value = 558440
$ws.Range("L131").Value = 1675320
$ws.Range("N131").Value = -1685400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 34263.75
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H70").Value = 9613.385
$ws.Range("I70").Value = 9165.166999999999
$ws.Range("K70").Value = 9165.166999999999
$ws.Range("M70").Value = -8895.166999999999
$ws.Range("H73").Value = 9613.385
$ws.Range("I73").Value = 9165.166999999999
$ws.Range("K73").Value = 9165.166999999999
$ws.Range("M73").Value = -8229.166999999999
$ws.Range("H80").Value = 1900
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 1900
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 922.86664
$ws.Range("I102").Value = 932
$ws.Range("J102").Value = 795
$ws.Range("K102").Value = 932
$ws.Range("L102").Value = 795
$ws.Range("M102").Value = 690
$ws.Range("N102").Value = -4039

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1887.3334
$ws.Range("I16").Value = 459.4
$ws.Range("K16").Value = 459.4
$ws.Range("M16").Value = -289.4
$ws.Range("H22").Value = 2753.1667
$ws.Range("I22").Value = 3353.8
$ws.Range("J22").Value = 2324.1428
$ws.Range("K22").Value = 3353.8
$ws.Range("L22").Value = 2324.1428
$ws.Range("M22").Value = -3058.8
$ws.Range("N22").Value = -2914.1428
$ws.Range("H27").Value = 2753.1667
$ws.Range("I27").Value = 3353.8
$ws.Range("J27").Value = 2324.1428
$ws.Range("K27").Value = 3353.8
$ws.Range("L27").Value = 2324.1428
$ws.Range("M27").Value = -3246.8
$ws.Range("N27").Value = -2538.1428
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 3172.182
$ws.Range("I40").Value = 1559
$ws.Range("K40").Value = 1559
$ws.Range("M40").Value = -1423
$ws.Range("H55").Value = 1739.8462
$ws.Range("I55").Value = 854
$ws.Range("J55").Value = 3157.2
$ws.Range("K55").Value = 854
$ws.Range("L55").Value = 3157.2
$ws.Range("M55").Value = -681
$ws.Range("N55").Value = -3503.2
$ws.Range("H63").Value = 90077
$ws.Range("I63").Value = 90077
$ws.Range("K63").Value = 90077
$ws.Range("M63").Value = -89328
$ws.Range("H66").Value = 90077
$ws.Range("I66").Value = 90077
$ws.Range("K66").Value = 270231
$ws.Range("M66").Value = -266487
$ws.Range("H68").Value = 2650.3333
$ws.Range("I68").Value = 2002
$ws.Range("J68").Value = 2780
$ws.Range("K68").Value = 2002
$ws.Range("L68").Value = 2780
$ws.Range("M68").Value = -1253
$ws.Range("N68").Value = -4278
$ws.Range("H71").Value = 2650.3333
$ws.Range("I71").Value = 2002
$ws.Range("J71").Value = 2780
$ws.Range("K71").Value = 10010
$ws.Range("L71").Value = 13900
$ws.Range("M71").Value = -6266
$ws.Range("N71").Value = -21388
$ws.Range("H74").Value = 79729
$ws.Range("I74").Value = 79729
$ws.Range("K74").Value = 79729
$ws.Range("M74").Value = -78731
$ws.Range("H77").Value = 79729
$ws.Range("I77").Value = 79729
$ws.Range("K77").Value = 239187
$ws.Range("M77").Value = -234195
$ws.Range("H82").Value = 1140.2
$ws.Range("I82").Value = 1175
$ws.Range("K82").Value = 1175
$ws.Range("M82").Value = -814
$ws.Range("H85").Value = 1140.2
$ws.Range("I85").Value = 1175
$ws.Range("K85").Value = 1175
$ws.Range("M85").Value = 73
$ws.Range("H100").Value = 4016.182
$ws.Range("I100").Value = 3459.875
$ws.Range("J100").Value = 5499.6665
$ws.Range("K100").Value = 3459.875
$ws.Range("L100").Value = 5499.6665
$ws.Range("M100").Value = -2918.875
$ws.Range("N100").Value = -6581.6665
$ws.Range("H122").Value = 6473.6553
$ws.Range("I122").Value = 6773.316
$ws.Range("J122").Value = 5904.3
$ws.Range("K122").Value = 20319.948
$ws.Range("L122").Value = 17712.9
$ws.Range("M122").Value = -17869.948
$ws.Range("N122").Value = -22612.9
$ws.Range("H132").Value = 2400.423
$ws.Range("I132").Value = 2144.85
$ws.Range("J132").Value = 3252.3333
$ws.Range("K132").Value = 6434.549999999999
$ws.Range("L132").Value = 9756.999899999999
$ws.Range("M132").Value = -3904.549999999999
$ws.Range("N132").Value = -14816.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1348
$ws.Range("H45").Value = 24128.334
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 24128.334
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 24128.334
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -25110.334
$ws.Range("H75").Value = 90118
$ws.Range("I75").Value = 90118
$ws.Range("K75").Value = 90118
$ws.Range("M75").Value = -89182
$ws.Range("H78").Value = 90118
$ws.Range("I78").Value = 90118
$ws.Range("K78").Value = 270354
$ws.Range("M78").Value = -265674
$ws.Range("H81").Value = 1667595.1
$ws.Range("I81").Value = 1113.8
$ws.Range("K81").Value = 2227.6
$ws.Range("M81").Value = -1166.6
$ws.Range("H84").Value = 1667595.1
$ws.Range("I84").Value = 1113.8
$ws.Range("K84").Value = 11138
$ws.Range("M84").Value = -5834
$ws.Range("H96").Value = 1577.8889
$ws.Range("I96").Value = 1333.6666
$ws.Range("J96").Value = 1700
$ws.Range("K96").Value = 1333.6666
$ws.Range("L96").Value = 1700
$ws.Range("M96").Value = 39.33339999999998
$ws.Range("N96").Value = -4446
$ws.Range("H100").Value = 6251515
$ws.Range("I100").Value = 6668082.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 13336165
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -13335624
$ws.Range("N100").Value = -7082
$ws.Range("H107").Value = 1239.3334
$ws.Range("J107").Value = 470
$ws.Range("L107").Value = 1410
$ws.Range("N107").Value = -5250
$ws.Range("H122").Value = 1040
$ws.Range("I122").Value = 961.25
$ws.Range("J122").Value = 1197.5
$ws.Range("K122").Value = 2883.75
$ws.Range("L122").Value = 3592.5
$ws.Range("M122").Value = -433.75
$ws.Range("N122").Value = -8492.5
$ws.Range("H126").Value = 6135
$ws.Range("I126").Value = 4736
$ws.Range("J126").Value = 7534
$ws.Range("K126").Value = 14208
$ws.Range("L126").Value = 22602
$ws.Range("M126").Value = -11738
$ws.Range("N126").Value = -27542
$ws.Range("H132").Value = 6458.25
$ws.Range("I132").Value = 7411.1113
$ws.Range("J132").Value = 3599.6667
$ws.Range("K132").Value = 22233.3339
$ws.Range("L132").Value = 10799.0001
$ws.Range("M132").Value = -19703.3339
$ws.Range("N132").Value = -15859.0001
$ws.Range("H136").Value = 4632.1177
$ws.Range("I136").Value = 4523.636
$ws.Range("J136").Value = 4831
$ws.Range("K136").Value = 13570.908
$ws.Range("L136").Value = 14493
$ws.Range("M136").Value = -11020.908
$ws.Range("N136").Value = -19593
